$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.671.20"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "3.511.98"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.36"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.79"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("D7").Value = "3.502.11"
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.196"
$ws.Range("E10").Value = "  +6.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.643"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.01"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000303"
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").Value = "4.071.14"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.29"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").Value = "69.661.53"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "3.517.63"
$ws.Range("E18").Value = "  -2.33%  "
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "533.68"
$ws.Range("E21").Value = "  +9.13%  "
$ws.Range("E22").Value = "  -3.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.20"
$ws.Range("E23").Value = "  -7.17%  "
$ws.Range("E24").Value = "  +5.49%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.56"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.09"
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.08"
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.11"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.40"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.90"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("E34").Value = "  -3.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "544.63"
$ws.Range("E35").Value = "  -5.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.407"
$ws.Range("E36").Value = "  +2.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.10"
$ws.Range("E37").Value = "  +5.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.14"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "0.0₃0760"
$ws.Range("E40").Value = "  -5.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.43"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("D43").Value = "3.353.77"
$ws.Range("E43").Value = "  +4.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.08"
$ws.Range("E44").Value = "  -4.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.49"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.96"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0439"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("E48").Value = "  -2.79%  "
$ws.Range("E49").Value = "  -7.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.79"
$ws.Range("E51").Value = "  +1.59%  "
